$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 2 de Abril de 2020 a las 23:50"

# Row 4
$ws.Cells.Item(4,2).Value = 240529
$ws.Cells.Item(4,3).Value = 25526
$ws.Cells.Item(4,5).Value = 224354

# Row 23
$ws.Cells.Item(23,1).Value = "Australia"
$ws.Cells.Item(23,2).Value = 5230
$ws.Cells.Item(23,3).Value = 182
$ws.Cells.Item(23,4).Value = 585
$ws.Cells.Item(23,5).Value = 4618
$ws.Cells.Item(23,6).Value = 50
$ws.Cells.Item(23,7).Value = 4
$ws.Cells.Item(23,8).Value = 27

# Row 24
$ws.Cells.Item(24,1).Value = "Noruega"
$ws.Cells.Item(24,2).Value = 5144
$ws.Cells.Item(24,3).Value = 267
$ws.Cells.Item(24,4).Value = 32
$ws.Cells.Item(24,5).Value = 5062
$ws.Cells.Item(24,6).Value = 96
$ws.Cells.Item(24,7).Value = 6
$ws.Cells.Item(24,8).Value = 50

# Row 37
$ws.Cells.Item(37,2).Value = 2421
$ws.Cells.Item(37,3).Value = 303
$ws.Cells.Item(37,4).Value = 125
$ws.Cells.Item(37,5).Value = 2262

# Row 71
$ws.Cells.Item(71,2).Value = 533
$ws.Cells.Item(71,3).Value = 74
$ws.Cells.Item(71,5).Value = 497

# Row 92
$ws.Cells.Item(92,4).Value = 76
$ws.Cells.Item(92,5).Value = 185

# Row 160
$ws.Cells.Item(160,1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(160,2).Value = 18
$ws.Cells.Item(160,3).Value = 2
$ws.Cells.Item(160,4).Value = 6
$ws.Cells.Item(160,5).Value = 11
$ws.Cells.Item(160,8).Value = 1

# Row 161
$ws.Cells.Item(161,1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(161,2).Value = 17
$ws.Cells.Item(161,4).Value = 0
$ws.Cells.Item(161,5).Value = 17

# Row 162
$ws.Cells.Item(162,1).Value = "Haiti"
$ws.Cells.Item(162,3).Value = 0
$ws.Cells.Item(162,4).Value = 1
$ws.Cells.Item(162,5).Value = 15
$ws.Cells.Item(162,8).Value = 0

# Row 163
$ws.Cells.Item(163,1).Value = "Siria"
$ws.Cells.Item(163,3).Value = 6
$ws.Cells.Item(163,4).Value = 0
$ws.Cells.Item(163,5).Value = 14
$ws.Cells.Item(163,8).Value = 2

# Row 166
$ws.Cells.Item(166,4).Value = 3
$ws.Cells.Item(166,5).Value = 11
